$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) from 2023-10-04 (45203)
# to 2023-10-06 (45205) for every data row (rows 2 through 360).
$ws.Range("C2:C360").Value = 45205
